# Apply the edits described by the upstream diff to the "Tabelle1" sheet
# of the workbook (the first sheet, which is the active/visible one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Make sure this sheet is the active one (it was already tabSelected in the
# original file) before changing the selected range.
$ws.Activate()

# --- Update cell values -------------------------------------------------

# Row 24: B24 17 -> 1
$ws.Range("B24").Value = 1

# Rows 34-39: B column 31 -> 32
$ws.Range("B34").Value = 32
$ws.Range("B35").Value = 32
$ws.Range("B36").Value = 32
$ws.Range("B37").Value = 32
$ws.Range("B38").Value = 32
$ws.Range("B39").Value = 32

# --- Update the current selection ---------------------------------------
# The sheet selection moves from A30:A46 (active cell A30) to the single
# cell B25 (active cell B25).
$ws.Range("B25").Select()
